# Updates the answer values in the two-digit ÷ one-digit division worksheet
# table. The document contains a single table; only the 5 non-empty rows
# (1, 5, 9, 13, 17) hold the "NN÷N=Q, R" answer strings, one per cell
# across 5 columns. Each cell's text is replaced in-place (preserving the
# existing run formatting) with its updated answer, per the commit diff.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "36÷7=5, 1"
$t.Cell(1,2).Range.Text  = "74÷2=37, 0"
$t.Cell(1,3).Range.Text  = "19÷8=2, 3"
$t.Cell(1,4).Range.Text  = "58÷6=9, 4"
$t.Cell(1,5).Range.Text  = "91÷2=45, 1"

$t.Cell(5,1).Range.Text  = "11÷9=1, 2"
$t.Cell(5,2).Range.Text  = "29÷5=5, 4"
$t.Cell(5,3).Range.Text  = "80÷4=20, 0"
$t.Cell(5,4).Range.Text  = "16÷4=4, 0"
$t.Cell(5,5).Range.Text  = "11÷3=3, 2"

$t.Cell(9,1).Range.Text  = "86÷5=17, 1"
$t.Cell(9,2).Range.Text  = "49÷4=12, 1"
$t.Cell(9,3).Range.Text  = "93÷5=18, 3"
$t.Cell(9,4).Range.Text  = "43÷3=14, 1"
$t.Cell(9,5).Range.Text  = "47÷4=11, 3"

$t.Cell(13,1).Range.Text = "77÷6=12, 5"
$t.Cell(13,2).Range.Text = "71÷4=17, 3"
$t.Cell(13,3).Range.Text = "96÷3=32, 0"
$t.Cell(13,4).Range.Text = "92÷2=46, 0"
$t.Cell(13,5).Range.Text = "22÷9=2, 4"

$t.Cell(17,1).Range.Text = "52÷8=6, 4"
$t.Cell(17,2).Range.Text = "71÷8=8, 7"
$t.Cell(17,3).Range.Text = "87÷2=43, 1"
$t.Cell(17,4).Range.Text = "45÷5=9, 0"
$t.Cell(17,5).Range.Text = "33÷8=4, 1"
